# Changed hardlookup for Node ID, to look for any column with name equal to ID
# Rename the "Name" header column to "ID" on the Server, Database, and Client sheets
# so the lookup can find the node-ID column by name.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Server")
$ws.Range("A1").Value = "ID"

$ws = $wb.Worksheets.Item("Database")
$ws.Range("A1").Value = "ID"

$ws = $wb.Worksheets.Item("Client")
$ws.Range("A1").Value = "ID"

# Server sheet: drop the two stray empty rows below the data (rows 4 and 6)
# and move the active-cell selection, matching the cleaned-up sheet saved
# after the header rename.
$ws = $wb.Worksheets.Item("Server")
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(4).Delete()
$ws.Range("B16").Select() | Out-Null
